# Updates the "RF001 - Autenticar Usuario" test-suite workbook from
# version 1.4 to 1.5: re-shuffles which expected-result / step text
# is associated with a handful of step rows across TC1, TC2, TC3 and TC4.
#
# Only the textual content of the affected cells is changed; everything
# else (styles, shared-string reuse, row/column layout) is left to Excel
# to manage automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# TC1 - step 2 (row 11): expected result changes from "TJSeg down" to
# "incorrect username/password"
$ws.Range("D11").Value = "SYSTEM alerta que o nome de usuario e/ou senha estao incorretos"

# TC1 - step 3 (row 12): step text changes from "select suggested user" to
# "fills in the fields"
$ws.Range("B12").Value = "Usuario do Sistema preenche os campos e clica no botao entrar"

# TC2 - step 2 (row 21): expected result changes from "CAS down" to
# "incorrect username/password"
$ws.Range("D21").Value = "SYSTEM alerta que o nome de usuario e/ou senha estao incorretos"

# TC2 - step 3 (row 22): step text changes from "fills in the fields" to
# "select suggested user"
$ws.Range("B22").Value = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"

# TC3 - step 2 (row 31): step text changes from "select suggested user" to
# "fills in the fields"; expected result changes from "incorrect
# username/password" to "TJSeg down"
$ws.Range("B31").Value = "Usuario do Sistema preenche os campos e clica no botao entrar"
$ws.Range("D31").Value = "SYSTEM alerta que o TJSeg (sistema que fornece as permissoes de acesso e escrita) esta fora do ar"

# TC3 - step 3 (row 32): step text changes from "fills in the fields" to
# "select suggested user"
$ws.Range("B32").Value = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"

# TC4 - step 2 (row 41): step text changes from "fills in the fields" to
# "select suggested user"; expected result changes from "incorrect
# username/password" to "CAS down"
$ws.Range("B41").Value = "Usuario do Sistema seleciona um nome de usuario sugerido, digita a senha e clica no botao entrar"
$ws.Range("D41").Value = "SYSTEM alerta que o CAS (sistema de autorizacao login-senha) esta fora do ar"

# TC6 - step 3 (row 62): step text changes from "select suggested user" to
# "fills in the fields"
$ws.Range("B62").Value = "Usuario do Sistema preenche os campos e clica no botao entrar"
